# Apply the scraped price/volume refresh for cryptos.xlsx (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 is an untouched, unstyled cell - copying its Style back onto each edited
# cell after the write keeps cells free of any extra "stored as text" styling
# (e.g. quote-prefix) that Excel would otherwise tack on for numeric-looking text.
$defaultStyle = $ws.Range("A1").Style

function Set-TextValue($cellRef, $text) {
    # Leading apostrophe forces Excel to store the value as text even when it
    # looks numeric (e.g. "586.68"), matching the source file's inlineStr cells.
    $ws.Range($cellRef).Value = "'" + $text
    $ws.Range($cellRef).Style = $defaultStyle
}

Set-TextValue 'D2' '69.637.11'
Set-TextValue 'E2' '  -1.03%  '
Set-TextValue 'D3' '3.509.91'
Set-TextValue 'E3' '  -2.51%  '
Set-TextValue 'E4' '  -0.04%  '
Set-TextValue 'D5' '586.68'
Set-TextValue 'E5' '  +1.11%  '
Set-TextValue 'D6' '183.32'
Set-TextValue 'E6' '  -3.59%  '
Set-TextValue 'D7' '3.497.57'
Set-TextValue 'E7' '  -2.73%  '
Set-TextValue 'E8' '  -3.13%  '
Set-TextValue 'D10' '0.197'
Set-TextValue 'E10' '  +6.33%  '
Set-TextValue 'E11' '  -2.84%  '
Set-TextValue 'D12' '53.98'
Set-TextValue 'E12' '  -3.81%  '
Set-TextValue 'E13' '  -2.28%  '
Set-TextValue 'D14' '9.45'
Set-TextValue 'E14' '  -2.61%  '
Set-TextValue 'D15' '4.068.19'
Set-TextValue 'D16' '19.27'
Set-TextValue 'E16' '  -2.74%  '
Set-TextValue 'D17' '69.594.13'
Set-TextValue 'E17' '  -1.00%  '
Set-TextValue 'D18' '3.510.51'
Set-TextValue 'E18' '  -2.47%  '
Set-TextValue 'E19' '  -2.66%  '
Set-TextValue 'D21' '531.59'
Set-TextValue 'E21' '  +8.06%  '
Set-TextValue 'E22' '  -3.62%  '
Set-TextValue 'D23' '18.22'
Set-TextValue 'E23' '  -7.06%  '
Set-TextValue 'D24' '4.58'
Set-TextValue 'E24' '  +4.82%  '
Set-TextValue 'E25' '  -0.90%  '
Set-TextValue 'D26' '95.45'
Set-TextValue 'E26' '  -1.12%  '
Set-TextValue 'D27' '11.11'
Set-TextValue 'E27' '  +0.82%  '
Set-TextValue 'E28' '  -0.95%  '
Set-TextValue 'D29' '9.09'
Set-TextValue 'E29' '  -3.18%  '
Set-TextValue 'D30' '32.18'
Set-TextValue 'E30' '  -0.23%  '
Set-TextValue 'D31' '7.29'
Set-TextValue 'E31' '  -4.18%  '
Set-TextValue 'D32' '12.41'
Set-TextValue 'E32' '  +1.14%  '
Set-TextValue 'D33' '63.85'
Set-TextValue 'E33' '  -3.56%  '
Set-TextValue 'E34' '  -3.90%  '
Set-TextValue 'D35' '544.33'
Set-TextValue 'E35' '  -5.65%  '
Set-TextValue 'D36' '3.13'
Set-TextValue 'E36' '  +6.27%  '
Set-TextValue 'E37' '  +2.10%  '
Set-TextValue 'D38' '38.02'
Set-TextValue 'E38' '  -2.02%  '
Set-TextValue 'E39' '  -0.12%  '
Set-TextValue 'D40' '0.0₃0759'
Set-TextValue 'E40' '  -6.89%  '
Set-TextValue 'E41' '  -2.24%  '
Set-TextValue 'D42' '3.39'
Set-TextValue 'E42' '  -2.78%  '
Set-TextValue 'D43' '3.355.22'
Set-TextValue 'E43' '  +4.05%  '
Set-TextValue 'E44' '  -5.02%  '
Set-TextValue 'B45' 'ApeXProtocol'
Set-TextValue 'C45' 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextValue 'D45' '3.50'
Set-TextValue 'E45' '  +2.90%  '
Set-TextValue 'B46' 'ThetaToken'
Set-TextValue 'C46' 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
Set-TextValue 'D46' '2.97'
Set-TextValue 'E46' '  -2.97%  '
Set-TextValue 'D47' '0.0439'
Set-TextValue 'E47' '  -1.83%  '
Set-TextValue 'E48' '  -3.06%  '
Set-TextValue 'D49' '8.96'
Set-TextValue 'E49' '  -7.16%  '
Set-TextValue 'E50' '  +0.03%  '
Set-TextValue 'D51' '137.77'
Set-TextValue 'E51' '  +2.74%  '
